# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worksheet lists, for each worker (Tipo Doc / N Doc / Nombre), the
# "Periodo Mora" (E) they owe for, plus Valor Mora (F) and Salario Basico (G).
# Previously the 6 workers were listed grouped by worker (period 1606 then
# 1605 for each). The sheet is reorganized to group by period instead
# (all 6 workers for period 1605, then all 6 workers for period 1606), and
# the Salario Basico values are refreshed to 689500 (the first worker had a
# stale 580000, and two rows for "DONELLA LOPEZ MARTINEZ" had a stale 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: TipoDoc, NroDoc, Nombre, PeriodoMora, ValorMora, SalarioBasico
$rows = @(
    @("CC", "1143348580", "KATHLEEN ELIANA TORRES LOPEZ",   "1605", 27600, 689500),
    @("CC", "1047436730", "SILVIA PATRICIA CASTILLO SAENZ", "1605", 27600, 689500),
    @("CC", "1081806502", "MARIA CAROLINA MENDOZA DIAZ",    "1605", 27600, 689500),
    @("CC", "1047438995", "VICTOR EDUARDO VASQUEZ ORTIZ",   "1605", 27600, 689500),
    @("CC", "45514581",   "DONELLA LOPEZ MARTINEZ",         "1605", 27600, 689500),
    @("CC", "1047444144", "RAFAEL ANDRES MOLINA GOMEZ",     "1605", 27600, 689500),
    @("CC", "1143348580", "KATHLEEN ELIANA TORRES LOPEZ",   "1606", 27600, 689500),
    @("CC", "1047436730", "SILVIA PATRICIA CASTILLO SAENZ", "1606", 27600, 689500),
    @("CC", "1081806502", "MARIA CAROLINA MENDOZA DIAZ",    "1606", 27600, 689500),
    @("CC", "1047438995", "VICTOR EDUARDO VASQUEZ ORTIZ",   "1606", 27600, 689500),
    @("CC", "45514581",   "DONELLA LOPEZ MARTINEZ",         "1606", 27600, 689500),
    @("CC", "1047444144", "RAFAEL ANDRES MOLINA GOMEZ",     "1606", 27600, 689500)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $data[0]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $data[1]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $data[2]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $data[3]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $data[4]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $data[5]   # G: Salario Basico
}
